# Support duplicate content by using title_hash - content_hash as TextID in excel file
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = "207af8cffc-15fc93b948"
    3  = "207af8cffc-45bc5d3397"
    4  = "207af8cffc-4c7e8364cf"
    5  = "207af8cffc-5c34b4bc84"
    6  = "207af8cffc-750f0cfa28"
    7  = "207af8cffc-753e0743b0"
    8  = "207af8cffc-89d8fe5c24"
    9  = "207af8cffc-aeda3f4cbe"
    10 = "207af8cffc-b5755d9969"
    11 = "207af8cffc-ff47c4299f"
    12 = "5682d36d17-02e8932823"
    13 = "5682d36d17-1002ece57f"
    14 = "5682d36d17-b5b9581ebe"
    15 = "5682d36d17-b735ff3ad3"
    16 = "5682d36d17-dc368eeb19"
}

foreach ($row in $values.Keys) {
    $val = $values[$row]
    $ws.Cells.Item($row, 1).Value = $val
    $ws.Cells.Item($row, 2).Value = $val
}
